$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue $ws.Cells.Item(2, 4) "27.649.11"
Set-TextValue $ws.Cells.Item(2, 5) "  +0.03%  "
Set-TextValue $ws.Cells.Item(3, 4) "1.845.06"
Set-TextValue $ws.Cells.Item(3, 5) "  -0.27%  "
Set-TextValue $ws.Cells.Item(4, 4) "1.011"
Set-TextValue $ws.Cells.Item(4, 5) "  -1.89%  "
Set-TextValue $ws.Cells.Item(5, 4) "317.41"
Set-TextValue $ws.Cells.Item(5, 5) "  -1.33%  "
Set-TextValue $ws.Cells.Item(6, 4) "1.009"
Set-TextValue $ws.Cells.Item(6, 5) "  -1.78%  "
Set-TextValue $ws.Cells.Item(7, 4) "0.4294"
Set-TextValue $ws.Cells.Item(7, 5) "  -1.90%  "
Set-TextValue $ws.Cells.Item(8, 4) "0.3743"
Set-TextValue $ws.Cells.Item(8, 5) "  -1.24%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.07317"
Set-TextValue $ws.Cells.Item(9, 5) "  -0.88%  "
Set-TextValue $ws.Cells.Item(10, 4) "0.8731"
Set-TextValue $ws.Cells.Item(10, 5) "  -1.06%  "
Set-TextValue $ws.Cells.Item(11, 4) "21.46"
Set-TextValue $ws.Cells.Item(11, 5) "  -0.31%  "
Set-TextValue $ws.Cells.Item(12, 4) "1.845.35"
Set-TextValue $ws.Cells.Item(12, 5) "  -0.28%  "
Set-TextValue $ws.Cells.Item(13, 4) "6.719"
Set-TextValue $ws.Cells.Item(13, 5) "  +0.02%  "
Set-TextValue $ws.Cells.Item(14, 4) "5.424"
Set-TextValue $ws.Cells.Item(14, 5) "  -1.34%  "
Set-TextValue $ws.Cells.Item(15, 4) "0.07102"
Set-TextValue $ws.Cells.Item(15, 5) "  -0.45%  "
Set-TextValue $ws.Cells.Item(16, 4) "88.63"
Set-TextValue $ws.Cells.Item(16, 5) "  +4.02%  "
Set-TextValue $ws.Cells.Item(17, 5) "  -1.83%  "
Set-TextValue $ws.Cells.Item(18, 4) "0.000008991"
Set-TextValue $ws.Cells.Item(18, 5) "  -0.88%  "
Set-TextValue $ws.Cells.Item(19, 4) "1.009"
Set-TextValue $ws.Cells.Item(19, 5) "  -1.76%  "
Set-TextValue $ws.Cells.Item(20, 4) "15.37"
Set-TextValue $ws.Cells.Item(20, 5) "  -0.48%  "
Set-TextValue $ws.Cells.Item(21, 4) "27.668.92"
Set-TextValue $ws.Cells.Item(21, 5) "  +0.04%  "
Set-TextValue $ws.Cells.Item(22, 4) "5.205"
Set-TextValue $ws.Cells.Item(22, 5) "  -1.61%  "
Set-TextValue $ws.Cells.Item(23, 4) "11.01"
Set-TextValue $ws.Cells.Item(23, 5) "  -2.43%  "
Set-TextValue $ws.Cells.Item(24, 4) "2.075.20"
Set-TextValue $ws.Cells.Item(24, 5) "  -0.32%  "
Set-TextValue $ws.Cells.Item(25, 4) "1.969"
Set-TextValue $ws.Cells.Item(25, 5) "  -2.97%  "
Set-TextValue $ws.Cells.Item(26, 4) "155.03"
Set-TextValue $ws.Cells.Item(26, 5) "  -1.70%  "
Set-TextValue $ws.Cells.Item(27, 4) "18.55"
Set-TextValue $ws.Cells.Item(27, 5) "  -0.79%  "
Set-TextValue $ws.Cells.Item(28, 4) "2.157"
Set-TextValue $ws.Cells.Item(28, 5) "  +8.02%  "
Set-TextValue $ws.Cells.Item(29, 4) "5.336"
Set-TextValue $ws.Cells.Item(29, 5) "  +0.09%  "
Set-TextValue $ws.Cells.Item(30, 4) "118.19"
Set-TextValue $ws.Cells.Item(30, 5) "  +0.21%  "
Set-TextValue $ws.Cells.Item(31, 4) "0.08910"
Set-TextValue $ws.Cells.Item(31, 5) "  -1.18%  "
Set-TextValue $ws.Cells.Item(32, 4) "1.221"
Set-TextValue $ws.Cells.Item(32, 5) "  +0.85%  "
Set-TextValue $ws.Cells.Item(33, 4) "0.7743"
Set-TextValue $ws.Cells.Item(33, 5) "  +0.34%  "
Set-TextValue $ws.Cells.Item(34, 4) "4.528"
Set-TextValue $ws.Cells.Item(34, 5) "  -0.47%  "
Set-TextValue $ws.Cells.Item(35, 4) "2.887"
Set-TextValue $ws.Cells.Item(35, 5) "  -3.41%  "
Set-TextValue $ws.Cells.Item(36, 4) "1.010"
Set-TextValue $ws.Cells.Item(36, 5) "  -1.78%  "
Set-TextValue $ws.Cells.Item(37, 4) "1.129"
Set-TextValue $ws.Cells.Item(37, 5) "  -1.12%  "
Set-TextValue $ws.Cells.Item(38, 4) "0.01979"
Set-TextValue $ws.Cells.Item(38, 5) "  +0.43%  "
Set-TextValue $ws.Cells.Item(39, 4) "0.05331"
Set-TextValue $ws.Cells.Item(39, 5) "  +0.96%  "
Set-TextValue $ws.Cells.Item(40, 4) "2.883"
Set-TextValue $ws.Cells.Item(40, 5) "  +1.42%  "
Set-TextValue $ws.Cells.Item(41, 4) "7.122"
Set-TextValue $ws.Cells.Item(41, 5) "  +4.15%  "
Set-TextValue $ws.Cells.Item(42, 4) "0.1691"
Set-TextValue $ws.Cells.Item(42, 5) "  +1.29%  "
Set-TextValue $ws.Cells.Item(43, 4) "0.5122"
Set-TextValue $ws.Cells.Item(43, 5) "  -1.12%  "
Set-TextValue $ws.Cells.Item(44, 4) "8.763"
Set-TextValue $ws.Cells.Item(44, 5) "  -0.16%  "
Set-TextValue $ws.Cells.Item(45, 4) "10.65"
Set-TextValue $ws.Cells.Item(45, 5) "  -0.54%  "
Set-TextValue $ws.Cells.Item(46, 4) "107.38"
Set-TextValue $ws.Cells.Item(46, 5) "  -2.55%  "
Set-TextValue $ws.Cells.Item(47, 4) "0.4744"
Set-TextValue $ws.Cells.Item(47, 5) "  +1.08%  "
Set-TextValue $ws.Cells.Item(48, 4) "0.06451"
Set-TextValue $ws.Cells.Item(48, 5) "  -2.28%  "
Set-TextValue $ws.Cells.Item(49, 4) "1.009"
Set-TextValue $ws.Cells.Item(49, 5) "  -1.96%  "
Set-TextValue $ws.Cells.Item(50, 4) "1.686"
Set-TextValue $ws.Cells.Item(50, 5) "  -1.02%  "
Set-TextValue $ws.Cells.Item(51, 4) "1.841"
Set-TextValue $ws.Cells.Item(51, 5) "  -2.74%  "
